# Turn the logs workbook into a single, read-only "Logs Sheet".
#
#  1. Drop the two unused blank worksheets (Sheet2 / Sheet3).
#  2. Rename the remaining Sheet1 -> "Logs Sheet" (the _FilterDatabase
#     defined name tracks the sheet and is updated automatically).
#  3. Protect the sheet (with a password) so the file behaves as read-only,
#     per the "Made EXCEL FILE read-only" commit message.

$wb = $excel.ActiveWorkbook

# 1. Remove the extra empty sheets.
[void]$wb.Worksheets.Item("Sheet2").Delete()
[void]$wb.Worksheets.Item("Sheet3").Delete()

# 2. Rename the data sheet.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Logs Sheet"

# 3. Protect the sheet (read-only) with a password.
$ws.Protect("D447")
